$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume data (Coin list refresh).
# Target cells originally hold plain text values (e.g. "1.001", "29.823.56",
# "  -1.42%  "), stored as inline strings.  Assigning a bare numeric-looking
# string via .Value would normally get auto-coerced into a real number by
# Excel, losing formatting like leading/trailing zeros or multi-dot
# thousands groups.  To keep them as genuine text (matching the original
# cell type) we force a Text number format before writing the value, then
# restore the default "Normal" style so no stray formatting is left behind.
$updates = @(
    @{ Cell = 'D2'; Value = '29.818.66' },
    @{ Cell = 'E2'; Value = '  -1.37%  ' },
    @{ Cell = 'D3'; Value = '1.893.32' },
    @{ Cell = 'E3'; Value = '  -1.03%  ' },
    @{ Cell = 'D4'; Value = '1.000' },
    @{ Cell = 'E4'; Value = '  -0.03%  ' },
    @{ Cell = 'D5'; Value = '0.7585' },
    @{ Cell = 'E5'; Value = '  +2.43%  ' },
    @{ Cell = 'D6'; Value = '239.34' },
    @{ Cell = 'E6'; Value = '  -1.66%  ' },
    @{ Cell = 'E7'; Value = '  +0.04%  ' },
    @{ Cell = 'D8'; Value = '1.891.77' },
    @{ Cell = 'E8'; Value = '  -0.41%  ' },
    @{ Cell = 'D9'; Value = '0.3043' },
    @{ Cell = 'E9'; Value = '  -3.22%  ' },
    @{ Cell = 'D10'; Value = '25.23' },
    @{ Cell = 'E10'; Value = '  -6.84%  ' },
    @{ Cell = 'D11'; Value = '0.06820' },
    @{ Cell = 'E11'; Value = '  -2.11%  ' },
    @{ Cell = 'D12'; Value = '0.07983' },
    @{ Cell = 'E12'; Value = '  +0.09%  ' },
    @{ Cell = 'D13'; Value = '0.7478' },
    @{ Cell = 'E13'; Value = '  -4.00%  ' },
    @{ Cell = 'D14'; Value = '1.885.15' },
    @{ Cell = 'E14'; Value = '  -2.17%  ' },
    @{ Cell = 'D15'; Value = '5.198' },
    @{ Cell = 'E15'; Value = '  -1.45%  ' },
    @{ Cell = 'D16'; Value = '91.11' },
    @{ Cell = 'E16'; Value = '  -0.50%  ' },
    @{ Cell = 'D17'; Value = '29.811.54' },
    @{ Cell = 'E17'; Value = '  -1.66%  ' },
    @{ Cell = 'D18'; Value = '5.986' },
    @{ Cell = 'E18'; Value = '  +3.15%  ' },
    @{ Cell = 'D19'; Value = '13.85' },
    @{ Cell = 'E19'; Value = '  -3.00%  ' },
    @{ Cell = 'D20'; Value = '0.000007660' },
    @{ Cell = 'E20'; Value = '  -2.09%  ' },
    @{ Cell = 'D21'; Value = '235.20' },
    @{ Cell = 'E21'; Value = '  -4.06%  ' },
    @{ Cell = 'E22'; Value = '  +0.11%  ' },
    @{ Cell = 'D23'; Value = '2.140.95' },
    @{ Cell = 'E23'; Value = '  -3.09%  ' },
    @{ Cell = 'D24'; Value = '1.000' },
    @{ Cell = 'E24'; Value = '  +0.04%  ' },
    @{ Cell = 'D25'; Value = '6.922' },
    @{ Cell = 'E25'; Value = '  +4.28%  ' },
    @{ Cell = 'D26'; Value = '9.246' },
    @{ Cell = 'E26'; Value = '  -1.64%  ' },
    @{ Cell = 'E27'; Value = '  +0.22%  ' },
    @{ Cell = 'D28'; Value = '18.68' },
    @{ Cell = 'E28'; Value = '  -1.76%  ' },
    @{ Cell = 'E29'; Value = '  +1.13%  ' },
    @{ Cell = 'D30'; Value = '2.056' },
    @{ Cell = 'E30'; Value = '  -3.55%  ' },
    @{ Cell = 'D31'; Value = '1.340' },
    @{ Cell = 'E31'; Value = '  -0.90%  ' },
    @{ Cell = 'D32'; Value = '1.511' },
    @{ Cell = 'E32'; Value = '  -2.28%  ' },
    @{ Cell = 'D33'; Value = '4.273' },
    @{ Cell = 'E33'; Value = '  -1.16%  ' },
    @{ Cell = 'D34'; Value = '4.016' },
    @{ Cell = 'D35'; Value = '0.05354' },
    @{ Cell = 'E35'; Value = '  +2.95%  ' },
    @{ Cell = 'D36'; Value = '1.244' },
    @{ Cell = 'E36'; Value = '  -4.20%  ' },
    @{ Cell = 'D37'; Value = '0.7277' },
    @{ Cell = 'E37'; Value = '  -3.33%  ' },
    @{ Cell = 'D38'; Value = '2.712' },
    @{ Cell = 'E38'; Value = '  -1.70%  ' },
    @{ Cell = 'D39'; Value = '0.01927' },
    @{ Cell = 'E39'; Value = '  -0.71%  ' },
    @{ Cell = 'E40'; Value = '  -0.87%  ' },
    @{ Cell = 'D41'; Value = '6.196' },
    @{ Cell = 'E41'; Value = '  -3.35%  ' },
    @{ Cell = 'D42'; Value = '0.4402' },
    @{ Cell = 'E42'; Value = '  -1.97%  ' },
    @{ Cell = 'D43'; Value = '72.32' },
    @{ Cell = 'E43'; Value = '  -4.86%  ' },
    @{ Cell = 'D44'; Value = '1.915' },
    @{ Cell = 'E44'; Value = '  -1.85%  ' },
    @{ Cell = 'D45'; Value = '1.001' },
    @{ Cell = 'E45'; Value = '  +0.13%  ' },
    @{ Cell = 'D46'; Value = '0.8228' },
    @{ Cell = 'E46'; Value = '  -1.30%  ' },
    @{ Cell = 'D47'; Value = '101.11' },
    @{ Cell = 'E47'; Value = '  -0.40%  ' },
    @{ Cell = 'D48'; Value = '7.582' },
    @{ Cell = 'E48'; Value = '  -1.16%  ' },
    @{ Cell = 'D49'; Value = '9.821' },
    @{ Cell = 'E49'; Value = '  -0.99%  ' },
    @{ Cell = 'D50'; Value = '2.045.64' },
    @{ Cell = 'E50'; Value = '  -3.52%  ' },
    @{ Cell = 'D51'; Value = '36.07' },
    @{ Cell = 'E51'; Value = '  -2.56%  ' }

)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
    $cell.Style = "Normal"
}
